$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 7 de Julio de 2020 a las 22:38"

# Update country data rows (refreshed case counts) and reordered country names
# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 3073971
$ws.Cells.Item(4, 3).Value = 33779
$ws.Cells.Item(4, 4).Value = 1341458
$ws.Cells.Item(4, 5).Value = 1598837
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 697
$ws.Cells.Item(4, 8).Value = 133676

# Row 8: Peru
$ws.Cells.Item(8, 1).Value = "Peru"
$ws.Cells.Item(8, 2).Value = 309278
$ws.Cells.Item(8, 3).Value = 3575
$ws.Cells.Item(8, 4).Value = 200938
$ws.Cells.Item(8, 5).Value = 97388
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = 180
$ws.Cells.Item(8, 8).Value = 10952

# Row 19: Alemania
$ws.Cells.Item(19, 1).Value = "Alemania"
$ws.Cells.Item(19, 2).Value = 198355
$ws.Cells.Item(19, 3).Value = 298
$ws.Cells.Item(19, 4).Value = 182700
$ws.Cells.Item(19, 5).Value = 6552
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 11
$ws.Cells.Item(19, 8).Value = 9103

# Row 27: Egipto
$ws.Cells.Item(27, 1).Value = "Egipto"
$ws.Cells.Item(27, 2).Value = 77279
$ws.Cells.Item(27, 3).Value = 1057
$ws.Cells.Item(27, 4).Value = 21718
$ws.Cells.Item(27, 5).Value = 52072
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(27, 7).Value = 67
$ws.Cells.Item(27, 8).Value = 3489

# Row 49: Israel
$ws.Cells.Item(49, 1).Value = "Israel"
$ws.Cells.Item(49, 2).Value = 32222
$ws.Cells.Item(49, 3).Value = 1473
$ws.Cells.Item(49, 4).Value = 18227
$ws.Cells.Item(49, 5).Value = 13653
$ws.Cells.Item(49, 6).Value = 0
$ws.Cells.Item(49, 7).Value = 8
$ws.Cells.Item(49, 8).Value = 342

# Row 62: Argelia
$ws.Cells.Item(62, 1).Value = "Argelia"
$ws.Cells.Item(62, 2).Value = 16879
$ws.Cells.Item(62, 3).Value = 475
$ws.Cells.Item(62, 4).Value = 12094
$ws.Cells.Item(62, 5).Value = 3817
$ws.Cells.Item(62, 6).Value = 0
$ws.Cells.Item(62, 7).Value = 9
$ws.Cells.Item(62, 8).Value = 968

# Row 70: Costa de Marfil
$ws.Cells.Item(70, 1).Value = "Costa de Marfil"
$ws.Cells.Item(70, 2).Value = 11194
$ws.Cells.Item(70, 3).Value = 228
$ws.Cells.Item(70, 4).Value = 5487
$ws.Cells.Item(70, 5).Value = 5631
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 1
$ws.Cells.Item(70, 8).Value = 76

# Row 89: Guinea
$ws.Cells.Item(89, 1).Value = "Guinea"
$ws.Cells.Item(89, 2).Value = 5636
$ws.Cells.Item(89, 3).Value = 26
$ws.Cells.Item(89, 4).Value = 4542
$ws.Cells.Item(89, 5).Value = 1060
$ws.Cells.Item(89, 6).Value = 0
$ws.Cells.Item(89, 7).Value = 0
$ws.Cells.Item(89, 8).Value = 34

# Row 90: Bosnia y Herzegovina
$ws.Cells.Item(90, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(90, 2).Value = 5621
$ws.Cells.Item(90, 3).Value = 163
$ws.Cells.Item(90, 4).Value = 2693
$ws.Cells.Item(90, 5).Value = 2721
$ws.Cells.Item(90, 6).Value = 0
$ws.Cells.Item(90, 7).Value = 8
$ws.Cells.Item(90, 8).Value = 207

# Row 98: Republica de Africa Central
$ws.Cells.Item(98, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(98, 2).Value = 4071
$ws.Cells.Item(98, 3).Value = 38
$ws.Cells.Item(98, 4).Value = 976
$ws.Cells.Item(98, 5).Value = 3043
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 52

# Row 106: Nicaragua
$ws.Cells.Item(106, 1).Value = "Nicaragua"
$ws.Cells.Item(106, 2).Value = 2846
$ws.Cells.Item(106, 3).Value = 327
$ws.Cells.Item(106, 4).Value = 1993
$ws.Cells.Item(106, 5).Value = 762
$ws.Cells.Item(106, 6).Value = 0
$ws.Cells.Item(106, 7).Value = 8
$ws.Cells.Item(106, 8).Value = 91

# Row 107: Mayotte
$ws.Cells.Item(107, 1).Value = "Mayotte"
$ws.Cells.Item(107, 2).Value = 2688
$ws.Cells.Item(107, 3).Value = 9
$ws.Cells.Item(107, 4).Value = 2446
$ws.Cells.Item(107, 5).Value = 208
$ws.Cells.Item(107, 6).Value = 0
$ws.Cells.Item(107, 7).Value = 0
$ws.Cells.Item(107, 8).Value = 34

# Row 111: Mali
$ws.Cells.Item(111, 1).Value = "Mali"
$ws.Cells.Item(111, 2).Value = 2348
$ws.Cells.Item(111, 3).Value = 17
$ws.Cells.Item(111, 4).Value = 1556
$ws.Cells.Item(111, 5).Value = 673
$ws.Cells.Item(111, 6).Value = 0
$ws.Cells.Item(111, 7).Value = 0
$ws.Cells.Item(111, 8).Value = 119

# Row 126: Cabo Verde
$ws.Cells.Item(126, 1).Value = "Cabo Verde"
$ws.Cells.Item(126, 2).Value = 1499
$ws.Cells.Item(126, 3).Value = 36
$ws.Cells.Item(126, 4).Value = 724
$ws.Cells.Item(126, 5).Value = 757
$ws.Cells.Item(126, 6).Value = 0
$ws.Cells.Item(126, 7).Value = 1
$ws.Cells.Item(126, 8).Value = 18

# Row 131: Ruanda
$ws.Cells.Item(131, 1).Value = "Ruanda"
$ws.Cells.Item(131, 2).Value = 1172
$ws.Cells.Item(131, 3).Value = 59
$ws.Cells.Item(131, 4).Value = 595
$ws.Cells.Item(131, 5).Value = 574
$ws.Cells.Item(131, 6).Value = 0
$ws.Cells.Item(131, 7).Value = 0
$ws.Cells.Item(131, 8).Value = 3

# Row 132: Jordania
$ws.Cells.Item(132, 1).Value = "Jordania"
$ws.Cells.Item(132, 2).Value = 1169
$ws.Cells.Item(132, 3).Value = 2
$ws.Cells.Item(132, 4).Value = 969
$ws.Cells.Item(132, 5).Value = 190
$ws.Cells.Item(132, 6).Value = 0
$ws.Cells.Item(132, 7).Value = 0
$ws.Cells.Item(132, 8).Value = 10

# Row 133: Letonia
$ws.Cells.Item(133, 1).Value = "Letonia"
$ws.Cells.Item(133, 2).Value = 1134
$ws.Cells.Item(133, 3).Value = 7
$ws.Cells.Item(133, 4).Value = 1008
$ws.Cells.Item(133, 5).Value = 96
$ws.Cells.Item(133, 6).Value = 0
$ws.Cells.Item(133, 7).Value = 0
$ws.Cells.Item(133, 8).Value = 30

# Row 134: Libia
$ws.Cells.Item(134, 1).Value = "Libia"
$ws.Cells.Item(134, 2).Value = 1117
$ws.Cells.Item(134, 3).Value = 0
$ws.Cells.Item(134, 4).Value = 269
$ws.Cells.Item(134, 5).Value = 814
$ws.Cells.Item(134, 6).Value = 0
$ws.Cells.Item(134, 7).Value = 0
$ws.Cells.Item(134, 8).Value = 34

# Row 135: Niger
$ws.Cells.Item(135, 1).Value = "Niger"
$ws.Cells.Item(135, 2).Value = 1094
$ws.Cells.Item(135, 3).Value = 1
$ws.Cells.Item(135, 4).Value = 974
$ws.Cells.Item(135, 5).Value = 52
$ws.Cells.Item(135, 6).Value = 0
$ws.Cells.Item(135, 7).Value = 0
$ws.Cells.Item(135, 8).Value = 68

# Row 137: Mozambique
$ws.Cells.Item(137, 1).Value = "Mozambique"
$ws.Cells.Item(137, 2).Value = 1040
$ws.Cells.Item(137, 3).Value = 28
$ws.Cells.Item(137, 4).Value = 280
$ws.Cells.Item(137, 5).Value = 752
$ws.Cells.Item(137, 6).Value = 0
$ws.Cells.Item(137, 7).Value = 0
$ws.Cells.Item(137, 8).Value = 8

# Row 145: Republica del Chad
$ws.Cells.Item(145, 1).Value = "Republica del Chad"
$ws.Cells.Item(145, 2).Value = 873
$ws.Cells.Item(145, 3).Value = 1
$ws.Cells.Item(145, 4).Value = 788
$ws.Cells.Item(145, 5).Value = 11
$ws.Cells.Item(145, 6).Value = 0
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 74

# Row 152: Togo
$ws.Cells.Item(152, 1).Value = "Togo"
$ws.Cells.Item(152, 2).Value = 689
$ws.Cells.Item(152, 3).Value = 9
$ws.Cells.Item(152, 4).Value = 467
$ws.Cells.Item(152, 5).Value = 207
$ws.Cells.Item(152, 6).Value = 0
$ws.Cells.Item(152, 7).Value = 0
$ws.Cells.Item(152, 8).Value = 15

# Row 209: Groenlandia
$ws.Cells.Item(209, 1).Value = "Groenlandia"
$ws.Cells.Item(209, 2).Value = 13
$ws.Cells.Item(209, 3).Value = 0
$ws.Cells.Item(209, 4).Value = 13
$ws.Cells.Item(209, 5).Value = 0
$ws.Cells.Item(209, 6).Value = 0
$ws.Cells.Item(209, 7).Value = 0
$ws.Cells.Item(209, 8).Value = 0

# Row 210: Islas Malvinas
$ws.Cells.Item(210, 1).Value = "Islas Malvinas"
$ws.Cells.Item(210, 2).Value = 13
$ws.Cells.Item(210, 3).Value = 0
$ws.Cells.Item(210, 4).Value = 13
$ws.Cells.Item(210, 5).Value = 0
$ws.Cells.Item(210, 6).Value = 0
$ws.Cells.Item(210, 7).Value = 0
$ws.Cells.Item(210, 8).Value = 0
